# This document embeds the same two logo pictures twice each: once in the
# "primary" header/footer and once in the "first page" header/footer of
# the section. Word had mislabelled the embedded picture parts, so the
# BTEC logo (alt text "BTec_Logo-Orange") is carrying the name
# "image2.jpg" and the Pearson logo (alt text ending in
# "...PearsonLogo.png") is carrying the name "image1.png". Fix the
# InlineShape names so they again match their actual picture parts:
#   BTec_Logo-Orange            : image2.jpg -> image1.jpg
#   ...PearsonLogo.png (descr)  : image1.png -> image2.png

$d = $word.ActiveDocument

function Fix-InlineShapeNames($range) {
    $shapes = $range.InlineShapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        $alt = $shape.AlternativeText

        if ($alt -eq "BTec_Logo-Orange") {
            try {
                $shape.Name = "image1.jpg"
            } catch {
                Write-Output ("Could not rename BTEC logo shape: " + $_.Exception.Message)
            }
        } elseif ($alt -like "*PearsonLogo.png") {
            try {
                $shape.Name = "image2.png"
            } catch {
                Write-Output ("Could not rename Pearson logo shape: " + $_.Exception.Message)
            }
        }
    }
}

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections.Item($s)

    # wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2,
    # wdHeaderFooterEvenPages = 3
    for ($hf = 1; $hf -le 3; $hf++) {
        $header = $section.Headers.Item($hf)
        if ($header.Exists) {
            Fix-InlineShapeNames $header.Range
        }

        $footer = $section.Footers.Item($hf)
        if ($footer.Exists) {
            Fix-InlineShapeNames $footer.Range
        }
    }
}
